# Applies the cryptos list update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.489.58"
$ws.Range("E2").Value = "  +1.62%  "
$ws.Range("D3").Value = "3.935.40"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "488.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.41%  "
$ws.Range("E7").Value = "  +1.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.997"
$ws.Range("D8").Style = "Normal"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000355"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.75%  "
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.34%  "
$ws.Range("D14").Value = "4.566.72"
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.75%  "
$ws.Range("D16").Value = "3.934.10"
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.136"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("E19").Value = "  -1.84%  "
$ws.Range("D20").Value = "68.563.25"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "445.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.68%  "
$ws.Range("E22").Value = "  +4.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +18.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.57%  "
$ws.Range("E27").Value = "  +3.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "720.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("E32").Value = "  -1.16%  "
$ws.Range("E33").Value = "  +2.94%  "
$ws.Range("D34").Value = "0.0₃0914"
$ws.Range("E34").Value = "  +14.32%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "42.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.50%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.20"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +15.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "60.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.150"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.403"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +19.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +15.62%  "
$ws.Range("E42").Value = "  +1.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.84%  "
$ws.Range("E45").Value = "  +1.05%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("E47").Value = "  +1.01%  "
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "146.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("D51").Value = "0.0₆0340"
$ws.Range("E51").Value = "  +35.03%  "
